# Weekly crime-data refresh for the 104th Precinct CompStat sheet.
# Bumps the report volume/week-range header text and overwrites the
# Crime Complaints table (rows 14-27, cols C:N) with the newly collected
# figures, including a few cells whose content flips between a literal
# "0"/"***.*" placeholder (shared text) and an actual number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text tweaks (rich-text cells) -- use Range.Replace so only the
# targeted cell is touched and the rest of the run's text is untouched.
# ---------------------------------------------------------------------
$ws.Range("A8").Replace("Number  6", "Number  7") | Out-Null
$ws.Range("C9").Replace("2/6/2023", "2/13/2023") | Out-Null
$ws.Range("C9").Replace("2/12/2023", "2/19/2023") | Out-Null

# ---------------------------------------------------------------------
# Helper: a handful of cells switch type between a text placeholder
# ("0" / "***.*", shared-string, General format) and a genuine number
# (#,##0 format). Plain ".Value = x" would keep the old number format,
# so first clone the destination cell's format from a same-styled
# neighbour via PasteSpecial(xlPasteFormats), then write the value
# (PasteSpecial(xlPasteValues) for text, so numeric-looking strings like
# "0" aren't auto-coerced back into a number).
# ---------------------------------------------------------------------
function Set-NumberCell($ref, $value) {
    $ws.Range("C16").Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($ref).Value = $value
}

function Set-TextPlaceholderCell($ref, $text) {
    $ws.Range("D14").Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null   # xlPasteFormats (General)
    if ($text -eq "0") {
        $ws.Range("D14").Copy() | Out-Null
    } else {
        $ws.Range("E14").Copy() | Out-Null
    }
    $ws.Range($ref).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# Row 15: C15 was the "0" placeholder, now a real count.
Set-NumberCell "C15" 2

# Row 22: C22 becomes numeric; D22/E22 become placeholders.
Set-NumberCell "C22" 1
Set-TextPlaceholderCell "D22" "0"
Set-TextPlaceholderCell "E22" "***.*"

# Row 26: C26 becomes numeric.
Set-NumberCell "C26" 2

# Row 27: C27/D27 become "0" placeholders, E27 becomes "***.*".
Set-TextPlaceholderCell "C27" "0"
Set-TextPlaceholderCell "D27" "0"
Set-TextPlaceholderCell "E27" "***.*"

# ---------------------------------------------------------------------
# Plain numeric overwrites (same type/format as before, value only).
# ---------------------------------------------------------------------
$ws.Range("N14").Value = -66.666666666666

$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = 50
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = -40

$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 31
$ws.Range("J16").Value = 24
$ws.Range("K16").Value = 29.166666666666
$ws.Range("L16").Value = 210
$ws.Range("M16").Value = -3.125
$ws.Range("N16").Value = -81.871345029239

$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 23
$ws.Range("H17").Value = 27.777777777777
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = 32
$ws.Range("K17").Value = 25
$ws.Range("L17").Value = 8.108108108108
$ws.Range("M17").Value = 48.148148148148
$ws.Range("N17").Value = 66.666666666666

$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -38.888888888888
$ws.Range("I18").Value = 27
$ws.Range("J18").Value = 36
$ws.Range("K18").Value = -25
$ws.Range("L18").Value = -10
$ws.Range("M18").Value = -55.737704918032
$ws.Range("N18").Value = -90.459363957597

$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -6.666666666666
$ws.Range("F19").Value = 62
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = 44.186046511627
$ws.Range("I19").Value = 86
$ws.Range("J19").Value = 83
$ws.Range("K19").Value = 3.614457831325
$ws.Range("L19").Value = 36.507936507936
$ws.Range("M19").Value = 45.762711864406
$ws.Range("N19").Value = 7.5

$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 118.181818181818
$ws.Range("I20").Value = 41
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = 41.379310344827
$ws.Range("L20").Value = 78.260869565217
$ws.Range("M20").Value = 24.242424242424
$ws.Range("N20").Value = -91.881188118811

$ws.Range("C21").Value = 46
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 43.75
$ws.Range("F21").Value = 143
$ws.Range("G21").Value = 107
$ws.Range("H21").Value = 33.644859813084
$ws.Range("I21").Value = 229
$ws.Range("J21").Value = 206
$ws.Range("K21").Value = 11.165048543689
$ws.Range("L21").Value = 40.490797546012
$ws.Range("M21").Value = 7.009345794392
$ws.Range("N21").Value = -78.618113912231

$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 3
$ws.Range("K22").Value = 50
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -25

$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 7.407407407407
$ws.Range("F24").Value = 85
$ws.Range("G24").Value = 98
$ws.Range("H24").Value = -13.265306122449
$ws.Range("I24").Value = 172
$ws.Range("J24").Value = 170
$ws.Range("K24").Value = 1.176470588235
$ws.Range("L24").Value = -14
$ws.Range("M24").Value = 10.256410256410

$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = -12.820512820512
$ws.Range("I25").Value = 61
$ws.Range("J25").Value = 61
$ws.Range("L25").Value = 3.389830508474
$ws.Range("M25").Value = -39

$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 3
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 200

$ws.Range("F27").Value = 7
$ws.Range("H27").Value = 40
